$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 262
$ws1.Range("F3").Value = 69
$ws1.Range("F5").Value = 6628
$ws1.Range("F6").Value = 5412
$ws1.Range("F7").Value = 449
$ws1.Range("F12").Value = 40

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 262
$ws4.Range("F3").Value = 69
$ws4.Range("F5").Value = 6628
$ws4.Range("F6").Value = 5412
$ws4.Range("F7").Value = 449
$ws4.Range("F14").Value = 40
